$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced into Text
# format before the assignment (otherwise Excel silently converts the string
# into a floating point number), and the style is then reset back to Normal
# so no extra formatting is left behind on the cell.
$ws.Range("D2").Value = "27.114.51"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "1.657.91"
$ws.Range("E3").Value = "  +3.79%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.892.43"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("D13").Value = "1.660.16"
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "27.100.04"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("D20").Value = "0.0₃0729"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.84%  "
$ws.Range("E23").Value = "  +4.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").Value = "1.523.74"
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("E35").Value = "  +9.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.889"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.32%  "
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.41%  "
$ws.Range("D44").Value = "1.798.60"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0977"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.06%  "
